$d = $word.ActiveDocument

# Update the title date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-07-23 Tuesday"

# Update all 100 equation cells in the table (20 rows x 5 columns), in document order
$values = @(
    "29+13=",
    "19-1=",
    "21+40=",
    "79-75=",
    "92-30=",
    "81-20=",
    "22+11=",
    "25+43=",
    "89-0=",
    "16+41=",
    "15+1=",
    "37-31=",
    "69-17=",
    "11+38=",
    "2+69=",
    "76+2=",
    "13+41=",
    "89-33=",
    "69-64=",
    "27-8=",
    "91-11=",
    "81-44=",
    "75-61=",
    "15+53=",
    "32-29=",
    "22-12=",
    "40-39=",
    "26+3=",
    "85-66=",
    "30+62=",
    "19+11=",
    "86+4=",
    "44-22=",
    "6+27=",
    "14+11=",
    "39-4=",
    "31+61=",
    "37+37=",
    "82-10=",
    "56+32=",
    "77-21=",
    "39+21=",
    "24+13=",
    "15+84=",
    "75+24=",
    "38+25=",
    "16+27=",
    "61-11=",
    "94-53=",
    "27+65=",
    "77-28=",
    "34-12=",
    "77-25=",
    "29+18=",
    "60-9=",
    "4+52=",
    "9+49=",
    "4+90=",
    "33+21=",
    "67-61=",
    "45+27=",
    "48+22=",
    "88+9=",
    "83-73=",
    "19+55=",
    "1+78=",
    "89-8=",
    "97-58=",
    "98-90=",
    "18+43=",
    "84-20=",
    "99-87=",
    "46-43=",
    "58+22=",
    "5+60=",
    "82-57=",
    "27+63=",
    "58-50=",
    "46-15=",
    "66-61=",
    "59-19=",
    "91-45=",
    "70-68=",
    "92-89=",
    "72-71=",
    "32+52=",
    "23+22=",
    "33+49=",
    "22+63=",
    "1+22=",
    "98-1=",
    "55-29=",
    "22+75=",
    "29+17=",
    "14+74=",
    "99-87=",
    "78-52=",
    "29-27=",
    "12+84=",
    "56-14="
)

$t = $d.Tables.Item(1)
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Output ("Updated cells: " + $idx)